# adding mock backend and change some layout
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window state: minimize the workbook window ---
$wb.Windows.Item(1).WindowState = -4140   # xlMinimized

# --- Row 4 (employee #2, Johnny Depay): nationality "england" -> "Australia" ---
$ws.Range("F4").Value = "Australia"

# --- New "suspend date" values ("-") for rows 5-14 (rows 3-4 already have it) ---
$ws.Range("K5").Value = "-"
$ws.Range("K6").Value = "-"
$ws.Range("K7").Value = "-"
$ws.Range("K8").Value = "-"
$ws.Range("K9").Value = "-"
$ws.Range("K10").Value = "-"
$ws.Range("K11").Value = "-"
$ws.Range("K12").Value = "-"
$ws.Range("K13").Value = "-"
$ws.Range("K14").Value = "-"

# --- New "office" column (P) ---
$ws.Range("P2").Value = "office"
$ws.Range("P3").Value = "Bali"
$ws.Range("P4").Value = "Jakarta"
$ws.Range("P5").Value = "Jogjakarta"
$ws.Range("P6").Value = "Jakarta"
$ws.Range("P7").Value = "Bali"
$ws.Range("P8").Value = "Jogjakarta"
$ws.Range("P9").Value = "Bandung"
$ws.Range("P10").Value = "Bandung"
$ws.Range("P11").Value = "Bandung"
$ws.Range("P12").Value = "Jogjakarta"
$ws.Range("P13").Value = "Jakarta"
$ws.Range("P14").Value = "Bali"

# --- Move the active selection to F3 ---
$ws.Range("F3").Select()
